# Allow generating multiple TeasuryData in single excel
#
# 1) Rename the sheet from "Treasury Example (with Math)" to the generic
#    name "sheet" so the template no longer hard-codes a single treasury
#    example name (callers can stamp out many of these).
# 2) Remove the "Par and Zero" line chart (and its drawing anchor) that
#    referenced the old sheet name / fixed data range, since it no longer
#    makes sense once the sheet can be duplicated/renamed per data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "sheet"

# Delete any chart objects (and their backing drawing anchors) on the sheet.
$charts = $ws.ChartObjects()
$chartCount = $charts.Count
for ($i = $chartCount; $i -ge 1; $i--) {
    $charts.Item($i).Delete()
}
